$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '31.136.13'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.78%  '
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.956.61'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '
$c.Style = "Normal"
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = "Normal"
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.32%  '
$c.Style = "Normal"
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '246.66'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +0.06%  '
$c.Style = "Normal"
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +0.19%  '
$c.Style = "Normal"
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4887'
$c.Style = "Normal"
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '
$c.Style = "Normal"
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2965'
$c.Style = "Normal"
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +1.71%  '
$c.Style = "Normal"
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06833'
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -2.20%  '
$c.Style = "Normal"
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '106.39'
$c.Style = "Normal"
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -5.04%  '
$c.Style = "Normal"
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.942.11'
$c.Style = "Normal"
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +0.23%  '
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +1.99%  '
$c.Style = "Normal"
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.409'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -1.19%  '
$c.Style = "Normal"
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.7148'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +4.80%  '
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '284.20'
$c.Style = "Normal"
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -5.02%  '
$c.Style = "Normal"
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '31.054.89'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +1.55%  '
$c.Style = "Normal"
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000007746'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +1.00%  '
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +0.72%  '
$c.Style = "Normal"
$c = $ws.Range('B20')
$c.NumberFormat = '@'
$c.Value = 'WrappedliquidstakedEther2.0'
$c.Style = "Normal"
$c = $ws.Range('C20')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '2.195.69'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +0.50%  '
$c.Style = "Normal"
$c = $ws.Range('B21')
$c.NumberFormat = '@'
$c.Value = 'Dai'
$c.Style = "Normal"
$c = $ws.Range('C21')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.Style = "Normal"
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +0.47%  '
$c.Style = "Normal"
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.524'
$c.Style = "Normal"
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -1.23%  '
$c.Style = "Normal"
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.42%  '
$c.Style = "Normal"
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.601'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +1.49%  '
$c.Style = "Normal"
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.916'
$c.Style = "Normal"
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +3.96%  '
$c.Style = "Normal"
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '168.94'
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.78%  '
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -2.77%  '
$c.Style = "Normal"
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.206'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +2.64%  '
$c.Style = "Normal"
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.55%  '
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -1.91%  '
$c.Style = "Normal"
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.733'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +15.23%  '
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.498'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +7.83%  '
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.04990'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -0.36%  '
$c.Style = "Normal"
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7621'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +2.05%  '
$c.Style = "Normal"
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +1.33%  '
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.02047'
$c.Style = "Normal"
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.19%  '
$c.Style = "Normal"
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.733'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.65%  '
$c.Style = "Normal"
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.719'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +0.96%  '
$c.Style = "Normal"
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '2.151'
$c.Style = "Normal"
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +5.44%  '
$c.Style = "Normal"
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '6.441'
$c.Style = "Normal"
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +10.05%  '
$c.Style = "Normal"
$c = $ws.Range('B41')
$c.NumberFormat = '@'
$c.Value = 'TrustWalletToken'
$c.Style = "Normal"
$c = $ws.Range('C41')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = "Normal"
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.8837'
$c.Style = "Normal"
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +1.18%  '
$c.Style = "Normal"
$c = $ws.Range('B42')
$c.NumberFormat = '@'
$c.Value = 'Quant'
$c.Style = "Normal"
$c = $ws.Range('C42')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = "Normal"
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '109.71'
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.59%  '
$c.Style = "Normal"
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.4459'
$c.Style = "Normal"
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -0.39%  '
$c.Style = "Normal"
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +3.89%  '
$c.Style = "Normal"
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.Style = "Normal"
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '7.509'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +2.63%  '
$c.Style = "Normal"
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '973.96'
$c.Style = "Normal"
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +14.96%  '
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +2.55%  '
$c.Style = "Normal"
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '9.376'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +0.75%  '
$c.Style = "Normal"
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.2609'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +2.26%  '
$c.Style = "Normal"
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +2.35%  '
$c.Style = "Normal"
